$d = $word.ActiveDocument

# --- 1a: merge the split runs (with gramStart/gramEnd proofErr markers around "A")
#     into a single run by doing a self Find/Replace over the whole sentence. A
#     Find.Execute replace collapses every run/proofErr inside the matched range
#     into one fresh run. ---
$text1a = "1a.  A man needs to get his cat, parrot and a bag of seed across a river using a boat that will only hold himself and one passenger."
$d.Content.Find.Execute($text1a, $true, $false, $false, $false, $false, $true, 1, $false, $text1a, 2) | Out-Null

# --- 1b: merge the split runs (with gramStart/gramEnd proofErr markers around "It")
#     into a single run the same way. ---
$text1b = "1b.  It doesn" + [char]8217 + "t say anything about how many trips it has to take and there is no way to only make 3 trips without leaving an incompatible pair together on one side or the other."
$d.Content.Find.Execute($text1b, $true, $false, $false, $false, $false, $true, 1, $false, $text1b, 2) | Out-Null

# --- 1c: add a brand-new response paragraph. The second-to-last paragraph in the
#     document is the (otherwise empty) paragraph that carries the _GoBack
#     bookmark. Insert a fresh empty paragraph right before it, then type the new
#     "1c." sentence into the start of that bookmark paragraph so the bookmark
#     stays attached to the end of the new text, just like in the target diff. ---
$bookmarkParaIndex = $d.Paragraphs.Count - 1
$bookmarkPara = $d.Paragraphs.Item($bookmarkParaIndex)

$rng = $bookmarkPara.Range
$rng.Collapse(1)
$rng.InsertParagraphBefore()

$bookmarkPara = $d.Paragraphs.Item($bookmarkParaIndex + 1)
$rng2 = $bookmarkPara.Range
$rng2.Collapse(1)
$text1c = "1c.  The overall goal is to get everything to the other side of the river and keep it all intact."
$rng2.InsertBefore($text1c)
